$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The "Metadata" sheet currently has 21 data rows (A1:B21). We need to insert a
# new "Jurisdiction" / "" row right after the existing "Contact" row (row 10),
# pushing the "Description" row (currently row 11) and everything below it down
# by one row, extending the sheet to A1:B22.

# 1) Pre-warm the formatting of the brand-new row 22 by copying the *format only*
#    from the current last row (21). Doing this before any full-content paste
#    keeps the destination cell's existing style index instead of minting a new,
#    unused one in styles.xml.
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Shift rows 11..21 down into rows 12..22, working from the bottom up so we
#    never overwrite a row before it has been copied out. The destination is
#    cleared first because PasteSpecial leaves a stale value behind in a
#    destination cell when the source cell being pasted is blank.
for ($r = 21; $r -ge 11; $r--) {
    $destRow = $r + 1
    $ws.Range("A" + $destRow + ":B" + $destRow).ClearContents()
    $ws.Range("A" + $r + ":B" + $r).Copy()
    $ws.Range("A" + $destRow + ":B" + $destRow).PasteSpecial(-4104)  # xlPasteAll
    $excel.CutCopyMode = $false
}

# 3) Write the new row 11 contents ("Jurisdiction" / empty value). The cell
#    already carries the correct body style (it held the old "Contact" row
#    before being overwritten in step 2), so no extra style work is needed.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# 4) Update the publication Date value (row 8, column B).
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"
